# feat: add 2022-Q1 data
#
# The existing "总计" (summary) sheet becomes the new "2022-Q1" sheet
# (reusing its sheetId/rId) and is repopulated with the per-fund holding
# data for 2022-Q1. A brand new "总计" sheet is then appended after it,
# containing the updated summary table (with the new 2022-Q1 row added
# on top of the previously existing quarters).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the current "总计" sheet into the "2022-Q1" sheet and
# fill it with the fund holdings table.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Header row (B1:H1) - reuse header style from an existing quarter sheet
$headerSrc = $wb.Worksheets.Item("2021-Q4")

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Data rows. Columns B,C,D,E,F,G are kept as text (matches source data,
# which stores these as strings, not numbers - e.g. "23.62", "010003",
# preserving formatting such as trailing zeros). Columns A (row index)
# and H (rank) are real numbers.
$q1Rows = @(
    @(0, "010003", "景顺长城电子信息产业股票A", "23.62", "93.46", "4.98", "1.1763", 7),
    @(1, "010792", "华安成长先锋混合A", "16.44", "93.43", "6.03", "0.9913", 5),
    @(2, "040015", "华安动态灵活配置混合", "22.93", "79.55", "3.15", "0.7223", 8),
    @(3, "006154", "华安制造先锋混合", "14.05", "93.81", "4.61", "0.6477", 6),
    @(4, "010004", "景顺长城电子信息产业股票C", "7.66", "93.46", "4.98", "0.3815", 7),
    @(5, "005823", "泰康颐享混合A", "14.39", "20.19", "1.70", "0.2446", 4),
    @(6, "010793", "华安成长先锋混合C", "3.80", "93.43", "6.03", "0.2291", 5),
    @(7, "001521", "国寿安保成长优选股票", "4.24", "87.83", "4.74", "0.2010", 4),
    @(8, "009490", "泰康科技创新一年定期开放混合", "2.61", "79.69", "7.17", "0.1871", 3),
    @(9, "159610", "景顺长城中证500增强策略ETF", "8.45", "98.35", "1.15", "0.0972", 9),
    @(10, "005824", "泰康颐享混合C", "2.82", "20.19", "1.70", "0.0479", 4),
    @(11, "008082", "国寿安保研究精选混合A", "0.52", "91.60", "3.62", "0.0188", 9),
    @(12, "008083", "国寿安保研究精选混合C", "0.15", "91.60", "3.62", "0.0054", 9),
    @(13, "006346", "安信量化优选股票A", "0.71", "90.62", "0.64", "0.0045", 6),
    @(14, "006347", "安信量化优选股票C", "0.49", "90.62", "0.64", "0.0031", 6)
)

foreach ($row in $q1Rows) {
    $r = [int]$row[0] + 2

    $q1.Range("B" + $r + ":G" + $r).NumberFormat = "@"

    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
}

# Apply the same style used elsewhere for the row-index column (A) to
# all the new data rows.
$headerSrc.Range("A2").Copy()
$q1.Range("A2:A16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: create a brand new "总计" sheet right after "2022-Q1" and
# populate it with the refreshed summary table.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalHeaderSrc = $wb.Worksheets.Item("2021-Q4")
$totalHeaderSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$totalRows = @(
    @(0, "2022-Q1", 15, 4.96),
    @(1, "2021-Q4", 40, 16.76),
    @(2, "2021-Q3", 54, 24.81),
    @(3, "2021-Q2", 63, 26.27),
    @(4, "2021-Q1", 22, 4.09),
    @(5, "2020-Q4", 24, 4.04)
)

foreach ($row in $totalRows) {
    $r = [int]$row[0] + 2
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

$totalHeaderSrc.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("A1").Select()
